# Apply updated dSF (column F) values as part of a data repull / mean
# recalculation pass. Only column F values change; everything else on
# the sheet stays the same.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Map of row number -> new value for column F ("dSF")
$updates = @{
    6  = 0
    8  = 0
    9  = -1
    15 = 0
    23 = -9
    27 = -1
    31 = 1
    32 = 0
    34 = 2
    35 = 3
    39 = 3
    42 = 2
    43 = 0
    48 = 2
    65 = -2
    67 = 0
    72 = 5
}

foreach ($row in $updates.Keys) {
    $ws.Range("F$row").Value = $updates[$row]
}
